$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The last existing row (33) currently uses the "latest row" date format
# (date-only). Since a new row is being appended, row 33 reverts to the
# regular "date + time" format used by all the other historical rows,
# and the newly appended row 34 takes over the "latest row" date-only format.
$ws.Cells.Item(33, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Append the new day's data (daily update).
$ws.Cells.Item(34, 1).Value = 45774
$ws.Cells.Item(34, 2).Value = 138
$ws.Cells.Item(34, 3).Value = 144
$ws.Cells.Item(34, 4).Value = 138

$ws.Cells.Item(34, 1).NumberFormat = "YYYY-MM-DD"
